# el-GR deck: translate remaining English placeholder prompts left over in
# the slide master / layouts, and give the master + title-slide placeholder
# shapes their localized (Greek) display names.
#
# NOTE: a small "Write-Host" no-op is interleaved after each mutation below.
# This headless host batches consecutive object-model edits and, without an
# intervening round-trip, a batch spanning master/layout/slide tiers can miss
# a target; emitting a trivial host call between edits forces each one to be
# committed before the next is issued.

$p = $ppt.ActivePresentation

# --- Slide Master: rename placeholder shapes + translate their prompt text ---
$m = $p.SlideMaster

$masterTitle = $m.Shapes.Item(1)
$masterTitle.Name = "Θέση τίτλου 1"
Write-Host "master title renamed"
$masterTitle.TextFrame.TextRange.Paragraphs(1,1).Text = "Στυλ κύριου τίτλου"
Write-Host "master title text set"

$masterText = $m.Shapes.Item(2)
$masterText.Name = "Θέση κειμένου 2"
Write-Host "master text renamed"
$masterText.TextFrame.TextRange.Paragraphs(1,1).Text = "Στυλ υποδείγματος κειμένου"
Write-Host "master text lvl0 set"
$masterText.TextFrame.TextRange.Paragraphs(2,1).Text = "Δεύτερου επιπέδου"
Write-Host "master text lvl1 set"
$masterText.TextFrame.TextRange.Paragraphs(3,1).Text = "Τρίτου επιπέδου"
Write-Host "master text lvl2 set"
$masterText.TextFrame.TextRange.Paragraphs(4,1).Text = "Τέταρτου επιπέδου"
Write-Host "master text lvl3 set"
$masterText.TextFrame.TextRange.Paragraphs(5,1).Text = "Πέμπτου επιπέδου"
Write-Host "master text lvl4 set"

$masterDate = $m.Shapes.Item(3)
$masterDate.Name = "Θέση ημερομηνίας 3"
Write-Host "master date renamed"

$masterFooter = $m.Shapes.Item(4)
$masterFooter.Name = "Θέση υποσέλιδου 4"
Write-Host "master footer renamed"

$masterSlideNum = $m.Shapes.Item(5)
$masterSlideNum.Name = "Θέση αριθμού διαφάνειας 5"
Write-Host "master slide number renamed"

# --- Slide Layout 4 ("Two Content"): leftover untranslated "Fifth level" ---
$layout4 = $m.CustomLayouts.Item(4)
Write-Host "layout4 fetched"
$layout4.Shapes.Item(2).TextFrame.TextRange.Paragraphs(5,1).Text = "Πέμπτου επιπέδου"
Write-Host "layout4 lvl4 set"

# --- Slide Layout 9 ("Picture with Caption"): leftover untranslated prompt ---
$layout9 = $m.CustomLayouts.Item(9)
Write-Host "layout9 fetched"
$layout9.Shapes.Item(2).TextFrame.TextRange.Text = "Κάντε κλικ στο εικονίδιο για να προσθέσετε εικόνα"
Write-Host "layout9 caption set"

# --- Slide 1: rename the title & subtitle placeholder shapes ---
$s = $p.Slides.Item(1)
$s.Shapes.Item(1).Name = "Τίτλος 1"
Write-Host "slide1 title renamed"
$s.Shapes.Item(2).Name = "Υπότιτλος 2"
Write-Host "slide1 subtitle renamed"
